$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data fixes (timeout-handling related corrections to the export) ---
$ws.Range("A2").Value = 1
$ws.Range("F2").Value = 9999

# --- New "Bug" flag column header in G1 ---
$ws.Range("G1").Value = "Bug"

# Copy the header formatting (bold font + centered/top alignment) from the
# neighbouring header cell F1 onto the new G1 header cell.
[void]$ws.Range("F1").Copy()
[void]$ws.Range("G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# G1 only needs left/right borders (no top/bottom), unlike the other headers,
# so drop the top and bottom edges inherited from F1's box border.
$topBorder = $ws.Range("G1").Borders.Item(8)
$topBorder.LineStyle = -4142
$bottomBorder = $ws.Range("G1").Borders.Item(9)
$bottomBorder.LineStyle = -4142

# Move the active selection to the new header cell.
[void]$ws.Range("G1").Select()
